$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "purpose" column (E2:E8): correct formatting from "fullRNASEQ" to "fullRNASeq"
for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    if ($cell.Text -ceq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
